# Updates the cryptocurrency price/volume table to reflect the latest
# Coinranking snapshot. Two coin pairs swapped rank order (Monero/Fetch.AI
# and Filecoin/Mantle and VeChain/RenzoRestakedETH), so those rows get
# their Coin/Link/Price/Volume columns rewritten; all other rows only get
# updated Price/Volume values.
#
# Note: several Price values look like plain numbers (e.g. "7.30",
# "529.46"). The source workbook stores them as text (inlineStr) so we
# prefix numeric-looking values with a leading apostrophe to force Excel
# to keep them as text instead of silently re-parsing/truncating them
# (e.g. "7.30" -> 7.3, "1.00" -> 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.428.62'
$ws.Range('E2').Value = '  -0.05%  '
$ws.Range('D3').Value = '3.100.62'
$ws.Range('E3').Value = '  -0.23%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''529.46'
$ws.Range('E5').Value = '  +0.88%  '
$ws.Range('D6').Value = '''137.59'
$ws.Range('E6').Value = '  +0.75%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = '3.101.12'
$ws.Range('E8').Value = '  -0.18%  '
$ws.Range('E9').Value = '  +4.54%  '
$ws.Range('D10').Value = '''7.30'
$ws.Range('E10').Value = '  +0.70%  '
$ws.Range('E11').Value = '  -0.08%  '
$ws.Range('D12').Value = '''0.411'
$ws.Range('E12').Value = '  +3.96%  '
$ws.Range('E13').Value = '  +1.36%  '
$ws.Range('D14').Value = '3.628.98'
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').Value = '''25.30'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = '''0.0000164'
$ws.Range('E16').Value = '  +0.26%  '
$ws.Range('D17').Value = '57.548.15'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').Value = '3.095.49'
$ws.Range('E18').Value = '  -0.34%  '
$ws.Range('D19').Value = '''5.99'
$ws.Range('E19').Value = '  +1.05%  '
$ws.Range('D20').Value = '''12.60'
$ws.Range('E20').Value = '  +1.15%  '
$ws.Range('D21').Value = '''8.05'
$ws.Range('E21').Value = '  +2.19%  '
$ws.Range('D22').Value = '''359.61'
$ws.Range('E22').Value = '  +3.49%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').Value = '''68.83'
$ws.Range('E24').Value = '  +1.03%  '
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = '''0.166'
$ws.Range('E26').Value = '  -0.83%  '
$ws.Range('E27').Value = '  -0.25%  '
$ws.Range('D28').Value = '0.0₃0862'
$ws.Range('E28').Value = '  -4.97%  '
$ws.Range('D29').Value = '''7.29'
$ws.Range('E29').Value = '  -1.74%  '
$ws.Range('E30').Value = '  -0.36%  '
$ws.Range('D31').Value = '''6.03'
$ws.Range('E31').Value = '  +0.81%  '
$ws.Range('D32').Value = '''21.31'
$ws.Range('E32').Value = '  +1.72%  '
$ws.Range('D33').Value = '''5.07'
$ws.Range('E33').Value = '  +2.38%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '''158.87'
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '''1.13'
$ws.Range('E35').Value = '  -2.30%  '
$ws.Range('D36').Value = '''6.03'
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('D37').Value = '''25.47'
$ws.Range('E37').Value = '  -1.72%  '
$ws.Range('E38').Value = '  +2.46%  '
$ws.Range('D39').Value = '''1.65'
$ws.Range('E39').Value = '  +2.20%  '
$ws.Range('D40').Value = '''0.0669'
$ws.Range('E40').Value = '  +0.96%  '
$ws.Range('D41').Value = '2.481.70'
$ws.Range('E41').Value = '  +5.67%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''3.98'
$ws.Range('E42').Value = '  -5.22%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '''0.694'
$ws.Range('E43').Value = '  -0.71%  '
$ws.Range('D44').Value = '''37.47'
$ws.Range('E44').Value = '  +2.74%  '
$ws.Range('D45').Value = '''1.00'
$ws.Range('E45').Value = '  +0.05%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '''0.0268'
$ws.Range('E46').Value = '  +0.19%  '
$ws.Range('B47').Value = 'RenzoRestakedETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D47').Value = '3.141.22'
$ws.Range('E47').Value = '  -0.11%  '
$ws.Range('D48').Value = '''0.983'
$ws.Range('E48').Value = '  +3.14%  '
$ws.Range('D49').Value = '''6.05'
$ws.Range('E49').Value = '  +0.40%  '
$ws.Range('D50').Value = '''19.70'
$ws.Range('E50').Value = '  -1.28%  '
$ws.Range('D51').Value = '''0.737'
$ws.Range('E51').Value = '  -3.06%  '
